$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate data rows 1849-1875 (A,B,C,D,J,K columns) ---
$ws.Range("A1848").Copy($ws.Range("A1849"))
$ws.Range("A1849").Value = 45423
$ws.Range("B1849").Value = "food"
$ws.Range("C1849").Value = "buttercroissant"
$ws.Range("D1849").Value = 0.59
$ws.Range("J1849").Value = "hofer"
$ws.Range("K1849").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1850"))
$ws.Range("A1850").Value = 45423
$ws.Range("B1850").Value = "food"
$ws.Range("C1850").Value = "laugenbrezel"
$ws.Range("D1850").Value = 0.45
$ws.Range("J1850").Value = "hofer"
$ws.Range("K1850").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1851"))
$ws.Range("A1851").Value = 45423
$ws.Range("B1851").Value = "food"
$ws.Range("C1851").Value = "hamburger mix 300g"
$ws.Range("D1851").Value = 1.49
$ws.Range("J1851").Value = "hofer"
$ws.Range("K1851").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1852"))
$ws.Range("A1852").Value = 45423
$ws.Range("B1852").Value = "food"
$ws.Range("C1852").Value = "blattsalat-mix"
$ws.Range("D1852").Value = 0.99
$ws.Range("J1852").Value = "hofer"
$ws.Range("K1852").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1853"))
$ws.Range("A1853").Value = 45422
$ws.Range("B1853").Value = "food"
$ws.Range("C1853").Value = "bio-datteltomaten"
$ws.Range("D1853").Value = 1.68
$ws.Range("J1853").Value = "hofer"
$ws.Range("K1853").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1854"))
$ws.Range("A1854").Value = 45422
$ws.Range("B1854").Value = "food"
$ws.Range("C1854").Value = "grie. Oliven gefu."
$ws.Range("D1854").Value = 1.69
$ws.Range("J1854").Value = "hofer"
$ws.Range("K1854").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1855"))
$ws.Range("A1855").Value = 45422
$ws.Range("B1855").Value = "food"
$ws.Range("C1855").Value = "buttercroissant"
$ws.Range("D1855").Value = 0.59
$ws.Range("J1855").Value = "hofer"
$ws.Range("K1855").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1856"))
$ws.Range("A1856").Value = 45422
$ws.Range("B1856").Value = "food"
$ws.Range("C1856").Value = "vegane vielfalt"
$ws.Range("D1856").Value = 1.86
$ws.Range("J1856").Value = "hofer"
$ws.Range("K1856").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1857"))
$ws.Range("A1857").Value = 45422
$ws.Range("B1857").Value = "food"
$ws.Range("C1857").Value = "salat gross"
$ws.Range("D1857").Value = 2.4
$ws.Range("J1857").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1858"))
$ws.Range("A1858").Value = 45423
$ws.Range("B1858").Value = "food"
$ws.Range("C1858").Value = "mill squeeze bbq"
$ws.Range("D1858").Formula = "=1.99/2"
$ws.Range("J1858").Value = "interspar"
$ws.Range("K1858").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1859"))
$ws.Range("A1859").Value = 45423
$ws.Range("B1859").Value = "food"
$ws.Range("C1859").Value = "sbudgernaturjo3,6%500g"
$ws.Range("D1859").Value = 0.79
$ws.Range("J1859").Value = "interspar"
$ws.Range("K1859").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1860"))
$ws.Range("A1860").Value = 45423
$ws.Range("B1860").Value = "food"
$ws.Range("C1860").Value = "sbudgernaturjo3,6%500g"
$ws.Range("D1860").Value = 0.79
$ws.Range("J1860").Value = "interspar"
$ws.Range("K1860").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1861"))
$ws.Range("A1861").Value = 45423
$ws.Range("B1861").Value = "food"
$ws.Range("C1861").Value = "bavaria blu"
$ws.Range("D1861").Formula = "=2.43/2"
$ws.Range("J1861").Value = "interspar"
$ws.Range("K1861").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1862"))
$ws.Range("A1862").Value = 45423
$ws.Range("B1862").Value = "food"
$ws.Range("C1862").Value = "spar highprotrote 500g"
$ws.Range("D1862").Value = 3.19
$ws.Range("J1862").Value = "interspar"
$ws.Range("K1862").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1863"))
$ws.Range("A1863").Value = 45425
$ws.Range("B1863").Value = "food"
$ws.Range("C1863").Value = "menu 2"
$ws.Range("D1863").Value = 2.9
$ws.Range("J1863").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1864"))
$ws.Range("A1864").Value = 45425
$ws.Range("B1864").Value = "food"
$ws.Range("C1864").Value = "sanlucar bananen (0,666kgx2.29eur/kg)"
$ws.Range("D1864").Value = 1.53
$ws.Range("J1864").Value = "billa"
$ws.Range("K1864").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1865"))
$ws.Range("A1865").Value = 45425
$ws.Range("B1865").Value = "food"
$ws.Range("C1865").Value = "lagencroissant"
$ws.Range("D1865").Value = 1.09
$ws.Range("J1865").Value = "billa"
$ws.Range("K1865").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1866"))
$ws.Range("A1866").Value = 45426
$ws.Range("B1866").Value = "food"
$ws.Range("C1866").Value = "menu 2"
$ws.Range("D1866").Value = 2.9
$ws.Range("J1866").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1867"))
$ws.Range("A1867").Value = 45426
$ws.Range("B1867").Value = "food"
$ws.Range("C1867").Value = "dessert"
$ws.Range("D1867").Value = 0.79
$ws.Range("J1867").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1868"))
$ws.Range("A1868").Value = 45427
$ws.Range("B1868").Value = "food"
$ws.Range("C1868").Value = "menu 2"
$ws.Range("D1868").Value = 2.9
$ws.Range("J1868").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1869"))
$ws.Range("A1869").Value = 45428
$ws.Range("B1869").Value = "food"
$ws.Range("C1869").Value = "menu 2"
$ws.Range("D1869").Value = 2.9
$ws.Range("J1869").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1870"))
$ws.Range("A1870").Value = 45428
$ws.Range("B1870").Value = "food"
$ws.Range("C1870").Value = "dessert"
$ws.Range("D1870").Value = 0.79
$ws.Range("J1870").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1871"))
$ws.Range("A1871").Value = 45429
$ws.Range("B1871").Value = "food"
$ws.Range("C1871").Value = "salat gross"
$ws.Range("D1871").Value = 2.4
$ws.Range("J1871").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1872"))
$ws.Range("A1872").Value = 45429
$ws.Range("B1872").Value = "food"
$ws.Range("C1872").Value = "dessert"
$ws.Range("D1872").Value = 0.79
$ws.Range("J1872").Value = "infineon"

$ws.Range("A1848").Copy($ws.Range("A1873"))
$ws.Range("A1873").Value = 45429
$ws.Range("B1873").Value = "food"
$ws.Range("C1873").Value = "apfel evelina (0,902kgx2.49eur/kg)"
$ws.Range("D1873").Formula = "=2.25-1.13"
$ws.Range("J1873").Value = "billa"
$ws.Range("K1873").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1874"))
$ws.Range("A1874").Value = 45429
$ws.Range("B1874").Value = "food"
$ws.Range("C1874").Value = "apfel evelina (0,902kgx2.49eur/kg)"
$ws.Range("D1874").Formula = "=2.25-1.13"
$ws.Range("J1874").Value = "billa"
$ws.Range("K1874").Value = "villach"

$ws.Range("A1848").Copy($ws.Range("A1875"))
$ws.Range("A1875").Value = 45429
$ws.Range("B1875").Value = "food"
$ws.Range("C1875").Value = "nuss schnecke"
$ws.Range("D1875").Value = 1.69
$ws.Range("J1875").Value = "billa"
$ws.Range("K1875").Value = "villach"

# --- Formula columns E,F,G,H,I ---
$ws.Range("E1849").Formula = "=MONTH(A1849)"
$ws.Range("F1849").Formula = "=YEAR(A1849)"
$ws.Range("G1849").Formula = "=WEEKDAY(A1849, 2)"
$ws.Range("H1849").Formula = "=CHOOSE(WEEKDAY(A1849, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1849").Formula = "=TEXT(A1849, ""MMM"")"

$ws.Range("E1850:E1863").Formula = "=MONTH(A1850)"
$ws.Range("F1850:F1863").Formula = "=YEAR(A1850)"
$ws.Range("G1850:G1863").Formula = "=WEEKDAY(A1850, 2)"
$ws.Range("H1850:H1863").Formula = "=CHOOSE(WEEKDAY(A1850, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1850:I1863").Formula = "=TEXT(A1850, ""MMM"")"

$ws.Range("E1864:E1866").Formula = "=MONTH(A1864)"
$ws.Range("F1864:F1866").Formula = "=YEAR(A1864)"
$ws.Range("G1864:G1866").Formula = "=WEEKDAY(A1864, 2)"
$ws.Range("H1864:H1866").Formula = "=CHOOSE(WEEKDAY(A1864, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1864:I1866").Formula = "=TEXT(A1864, ""MMM"")"

$ws.Range("E1867:E1869").Formula = "=MONTH(A1867)"
$ws.Range("F1867:F1869").Formula = "=YEAR(A1867)"
$ws.Range("G1867:G1869").Formula = "=WEEKDAY(A1867, 2)"
$ws.Range("H1867:H1869").Formula = "=CHOOSE(WEEKDAY(A1867, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1867:I1869").Formula = "=TEXT(A1867, ""MMM"")"

$ws.Range("E1870:E1871").Formula = "=MONTH(A1870)"
$ws.Range("F1870:F1871").Formula = "=YEAR(A1870)"
$ws.Range("G1870:G1871").Formula = "=WEEKDAY(A1870, 2)"
$ws.Range("H1870:H1871").Formula = "=CHOOSE(WEEKDAY(A1870, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1870:I1871").Formula = "=TEXT(A1870, ""MMM"")"

$ws.Range("E1872").Formula = "=MONTH(A1872)"
$ws.Range("F1872").Formula = "=YEAR(A1872)"
$ws.Range("G1872").Formula = "=WEEKDAY(A1872, 2)"
$ws.Range("H1872").Formula = "=CHOOSE(WEEKDAY(A1872, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1872").Formula = "=TEXT(A1872, ""MMM"")"

$ws.Range("E1873:E1875").Formula = "=MONTH(A1873)"
$ws.Range("F1873:F1875").Formula = "=YEAR(A1873)"
$ws.Range("G1873:G1875").Formula = "=WEEKDAY(A1873, 2)"
$ws.Range("H1873:H1875").Formula = "=CHOOSE(WEEKDAY(A1873, 2), ""Monday"", ""Tuesday"",""Wednesday"", ""Thursday"", ""Friday"", ""Saturday"",""Sunday"")"
$ws.Range("I1873:I1875").Formula = "=TEXT(A1873, ""MMM"")"

# --- View state: update frozen-pane scroll position and selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1855
$win.ScrollColumn = 1
$ws.Range("I1870").Select()
